$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column D (sex = "U") for rows 2-11
$ws.Range("D2:D11").Value = "U"

# Fill column F (sire) and G (dam) values
$sireVals = @(1, 3, 5, 3, 5, 6, 8, 10, 8, 10)
$damVals  = @(1, 1, 1, 5, 5, 6, 6, 6, 10, 10)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $sireVals[$i]
    $ws.Cells.Item($row, 7).Value = $damVals[$i]
}

# Update the selection to match the diff
$ws.Range("D2:G11").Select()
